$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily observations appended to the bottom of the series (MV -datos- update)
$data = @(
    @("09-09-2021", 1.48, 2.03, 2.54, 3.16, -0.46),
    @("10-09-2021", 1.5,  2.11, 2.59, 3.42, -0.35),
    @("13-09-2021", 1.53, 2.23, 2.69, 3.61, -0.34),
    @("14-09-2021", 1.56, 2.24, 2.74, 3.59, -0.34)
)

$startRow = 175
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]

    # Column A holds the date label as text (e.g. "09-09-2021"), not an Excel date
    # serial, matching the rest of the "Serie" column. Force text type via the
    # number format, then restore the default (unstyled) cell style so the new
    # row's formatting matches the other data rows.
    $dateCell = $ws.Cells.Item($row, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $values[0]
    $dateCell.Style = "Normal"

    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
}
